$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert the "Featured" column (S) for rows 5-9 from "Y" back to "N"
$ws.Range("S5:S9").Value = "N"

# Restore the selection/active cell to J12 (and implicitly clear the
# scrolled topLeftCell state left over from the previous selection at S16)
$ws.Range("J12").Select()
